# Heston SL results workbook: new power-meter (laptop) measurement block +
# date/version comparison against the earlier run. See commit message:
# "new pwer meter laptop & date comparison"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$GRAY = 8355711   # RGB(127,127,127) ~ Text1, Lighter 50% theme tint used by the old block

# ---------------------------------------------------------------------------
# 1. Wipe the old "second half" of the sheet (rows 21-40) - we rebuild it at
#    its new, shifted locations below.
# ---------------------------------------------------------------------------
$ws.Range("A21:E40").Clear()

# ---------------------------------------------------------------------------
# 2. Old measurement block (originally rows 22-33) now lives one row lower,
#    rows 23-33, and is re-styled in gray to mark it as the "old" run.
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "Performance (RunCPU):"
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").Font.Color = $GRAY

$ws.Range("A24").Value = "Version: 052cf806968faa45a7c5b80d785a06e9640806fd"
$ws.Range("A24").Font.Color = $GRAY

$ws.Range("A25").Value = 557568
$ws.Range("A25").Font.Color = $GRAY
$ws.Range("B25").Value = "values / sec"
$ws.Range("B25").Font.Color = $GRAY

$ws.Range("A26").Formula = "=A25*250"
$ws.Range("A26").NumberFormat = "##0.0E+0"
$ws.Range("A26").Font.Color = $GRAY
$ws.Range("B26").Value = "steps / sec"
$ws.Range("B26").Font.Color = $GRAY

$ws.Range("A27").NumberFormat = "0.00E+00"
$ws.Range("A27").Font.Color = $GRAY

$ws.Range("A28").Value = "Power (Einfaches Wattmeter, RunCPU):"
$ws.Range("A28").Font.Bold = $true
$ws.Range("A28").Font.Color = $GRAY

$ws.Range("A29").Value = 35.5
$ws.Range("A29").Font.Color = $GRAY
$ws.Range("B29").Value = "Watt"
$ws.Range("B29").Font.Color = $GRAY

$ws.Range("A30").Value = "No display"
$ws.Range("A30").Font.Color = $GRAY

$ws.Range("A32").Value = "Power Efficiency (RunCPU):"
$ws.Range("A32").Font.Bold = $true
$ws.Range("A32").Font.Color = $GRAY

$ws.Range("A33").Formula = "=A29/A26"
$ws.Range("A33").NumberFormat = "##0.0E+0"
$ws.Range("A33").Font.Color = $GRAY
$ws.Range("B33").Value = "J / step"
$ws.Range("B33").Font.Color = $GRAY

# ---------------------------------------------------------------------------
# 3. New block: latest run performance + version, on today's laptop.
#    (the version cell is the one that used to read "Version: " - update it
#    first so it keeps its place in the shared string table, matching the
#    order these strings were authored in)
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = "Version: 22425123e417c44fe57514d14b7c0cf6c4ca185b"

$ws.Range("A36").Value = "Performance (run_cpu):"
$ws.Range("A36").Font.Bold = $true

$ws.Range("A38").Value = 141086000
$ws.Range("A38").NumberFormat = "##0.0E+0"
$ws.Range("A38").HorizontalAlignment = -4152
$ws.Range("B38").Value = "steps / sec"

# ---------------------------------------------------------------------------
# 4. Power (Volcraft VC 870) block, shifted down, plus the new laptop
#    measurement columns (D/E) on the run_cpu row.
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = "Power (Volcraft VC 870, run_cpu)"
$ws.Range("A40").Font.Bold = $true

$ws.Range("A41").Value = "Nur Netzteil"
$ws.Range("C41").Value = 0

$ws.Range("A42").Value = "Idle, display aus"
$ws.Range("C42").Value = 12.9

$ws.Range("A43").Value = "Idle, display aus, Batterie raus, Netwerk aus, Wifi aus"
$ws.Range("C43").Value = 11.2

$ws.Range("A44").Value = "run_cpu"
$ws.Range("C44").Value = 30.6
$ws.Range("D44").Value = "Watt"
$ws.Range("E44").Value = "Display aus, Batterie raus, Netwerk aus, Wifi aus, Reboot & Warmup, Multiple starts"

# ---------------------------------------------------------------------------
# 5. New Power Efficiency section for the Volcraft measurement.
# ---------------------------------------------------------------------------
$ws.Range("A46").Value = "Power Efficiency (RunCPU):"
$ws.Range("A46").Font.Bold = $true

$ws.Range("A47").Formula = "=C44/A38"
$ws.Range("A47").NumberFormat = "##0.0E+0"
$ws.Range("B47").Value = "J / step"

# ---------------------------------------------------------------------------
# 6. View state: selection like the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("E45").Select()

Write-Host "edit applied"
